# Burndown chart para el reporte de avance del 01/11/2010
#
# - Mark all remaining (not-yet-completed) Sprint tasks as "completed"
#   in the Sprint sheet (column C, rows 7-14), which also causes the
#   now-unused "active"/"pending" shared strings to be dropped on save.
# - Make "Sprint" the active sheet/tab again (it was "Burndown Chart"),
#   and update the selected cell on each sheet accordingly.

$wb = $excel.ActiveWorkbook

$sprint = $wb.Worksheets.Item("Sprint")

# Update task statuses: everything that was "active" or "pending" is now
# "completed". The "Burndown Chart" sheet's selection (B24) is left as-is.
$sprint.Range("C7").Value = "completed"
$sprint.Range("C8").Value = "completed"
$sprint.Range("C9").Value = "completed"
$sprint.Range("C10").Value = "completed"
$sprint.Range("C11").Value = "completed"
$sprint.Range("C12").Value = "completed"
$sprint.Range("C13").Value = "completed"
$sprint.Range("C14").Value = "completed"

# Move the active tab/selection from "Burndown Chart" back to "Sprint".
$sprint.Activate()
$sprint.Range("C15").Select()
